# "Variables categóricas a one hot vectors ..." bullet gets struck through
# (marked as done/obsolete) and its trailing "Ej: ..." example sentence,
# which used to be interrupted mid-word by a stray "_GoBack" bookmark, is
# tidied up so the bookmark now wraps the whole struck-through sentence
# instead of splitting "consideradas" into "con" + "sideradas".

$d = $word.ActiveDocument

# --- Step 1: heal the run the old bookmark used to split --------------
# In the original file <w:bookmarkStart w:name="_GoBack"/> sits between
# "...posibilidades con" and "sideradas en la variable. ", i.e. in the
# middle of the word "consideradas". Re-"typing" that exact text via
# Find/Replace collapses it back into a single contiguous run and leaves
# the (now zero-width) bookmark sitting right after it.
$healRange = $d.Content
$healRange.Find.Execute("consideradas en la variable. ", $false, $false, $false, $false, $false, $true, 1, $false, "consideradas en la variable. ", 2) | Out-Null

# --- Step 2: drop the old _GoBack bookmark -----------------------------
# It will be re-created below spanning the whole sentence that gets
# struck through.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 3: strike through the whole (now obsolete) explanation -------
$target = $d.Content
$found = $target.Find.Execute("Variables categóricas a one hot vectors, cuyo tamaño sea igual al número de posibilidades consideradas en la variable. Ej: meses del año, vector de 12 vars que pueden ser o uno o cero", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $target.Font.StrikeThrough = 1

    # --- Step 4: re-plant _GoBack around the struck-through sentence ---
    $d.Bookmarks.Add("_GoBack", $target)
}

Write-Output "one-hot bullet struck through: $found"
